# "Luận thêm phần đối của cung vô chính diệu"
# Adds the "tại cung đối Huynh Đệ" ("at the opposite palace of Huynh Đệ")
# readings for the 14 main stars (solo "tọa thủ") plus every pairwise
# "đồng cung" combination, appended after the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$stars = @("Tử Vi", "Thiên Cơ", "Thái Dương", "Vũ Khúc", "Thiên Đồng", "Liêm Trinh", "Thiên Phủ", "Thái Âm", "Tham Lang", "Cự Môn", "Thiên Tướng", "Thiên Lương", "Thất Sát", "Phá Quân")

# New data starts at row 4291 (row 4290 is intentionally left blank/untouched).
$row = 4291

for ($i = 0; $i -lt $stars.Length; $i++) {
    $text = $stars[$i] + " tọa thủ tại cung đối Huynh Đệ"
    $ws.Cells.Item($row, 1).Value = $text
    $ws.Cells.Item($row, 2).Value = $text
    $row = $row + 1
}

for ($i = 0; $i -lt $stars.Length; $i++) {
    for ($j = $i + 1; $j -lt $stars.Length; $j++) {
        $text = $stars[$i] + " đồng cung " + $stars[$j] + " tại cung đối Huynh Đệ"
        $ws.Cells.Item($row, 1).Value = $text
        $ws.Cells.Item($row, 2).Value = $text
        $row = $row + 1
    }
}

$ws.Range("I4312").Select()
